# Updated symbol list on Thu Feb 16 18:55:32 UTC 2023 with GitHub Actions
#
# This refreshes the "Price" (column D) and "Volume(1h)" (column E) readings
# for the coin rows that moved since the last scrape. The sheet stores these
# numbers/percentages as literal text (inline strings) rather than numeric
# cells, so each target cell is forced to Text format ("@") before the new
# value is written -- otherwise Excel would auto-convert a value like
# "6.16%" into the number 0.0616 with a percentage format, which is not what
# the source data looks like.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "321.51"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "6.16%"

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "49.31"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "11.58%"

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.305"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "3.86%"

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.08064"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "4.26%"

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "4.585"

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.346"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "28.45%"

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.653"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "2.31%"

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.1267"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-1.23%"

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1971"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "5.68%"

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.09684"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "4.45%"

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.04696"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "13.13%"

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.1049"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "0.00%"

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.001322"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "2.96%"

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.04205"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "0.37%"

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.005902"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "2.56%"

# Row 17
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-0.08%"

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.498"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "7.19%"

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.3522"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "5.08%"

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "8.143"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "1.29%"

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.1380"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "0.31%"

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.3088"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "-2.83%"

# Row 23
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "1.28%"

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.004292"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-2.54%"

# Row 25
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "0.12%"

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0003536"

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02724"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "8.59%"

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05974"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "12.55%"

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01078"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "85.06%"

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.008034"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "3.91%"

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1470"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "8.61%"

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.007758"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "5.49%"

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.007876"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "4.92%"

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.3214"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "6.50%"

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006949"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "4.12%"

# Row 47
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "0.14%"

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.05594"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "-13.01%"

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.003996"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "-4.87%"

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002098"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.14%"

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0001998"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "0.14%"
